$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of configuration data
$ws.Range("A3").Value = "URL"
$ws.Range("B3").Value = "https://itdashboard.gov/"
$ws.Range("A4").Value = "Directory Name"
$ws.Range("B4").Value = "output"

# C4 gets the same (underline) font formatting as E6
$ws.Range("C4").Font.Underline = $true

# Widen column A to fit the new, longer labels
$ws.Columns.Item(1).ColumnWidth = 14

# Update selection to C4
$ws.Range("C4").Select()
